$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.312.14"
$ws.Range("E2").Value = "  +0.35%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.668.34"
$ws.Range("E3").Value = "  +3.16%  "

# Row 4
$ws.Range("E4").Value = "  -0.02%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "608.81"
$ws.Range("E5").Value = "  +4.36%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "143.60"
$ws.Range("E6").Value = "  -1.03%  "

# Row 7
$ws.Range("E7").Value = "  -0.01%  "

# Row 8
$ws.Range("E8").Value = "  -0.98%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.668.16"
$ws.Range("E9").Value = "  +3.16%  "

# Row 10
$ws.Range("E10").Value = "  +0.15%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.64"
$ws.Range("E11").Value = "  +0.74%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.153"
$ws.Range("E12").Value = "  +0.69%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.358"
$ws.Range("E13").Value = "  +1.86%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.40"
$ws.Range("E14").Value = "  +1.03%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.146.16"
$ws.Range("E15").Value = "  +3.04%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "63.131.00"
$ws.Range("E16").Value = "  +0.17%  "

# Row 17
$ws.Range("E17").Value = "  -0.35%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.666.81"
$ws.Range("E18").Value = "  +3.46%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.44"
$ws.Range("E19").Value = "  +2.62%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "340.30"
$ws.Range("E20").Value = "  -0.64%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.40"
$ws.Range("E21").Value = "  +1.09%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.86"
$ws.Range("E22").Value = "  +3.07%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "67.57"
$ws.Range("E24").Value = "  +0.05%  "

# Row 25
$ws.Range("E25").Value = "  +2.53%  "

# Row 26
$ws.Range("E26").Value = "  -2.66%  "

# Row 27
$ws.Range("E27").Value = "  -0.18%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.58"
$ws.Range("E28").Value = "  +3.58%  "

# Row 29
$ws.Range("B29").Value = "Bittensor"
$ws.Range("C29").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "541.36"
$ws.Range("E29").Value = "  +15.69%  "

# Row 30
$ws.Range("B30").Value = "Binance-PegBSC-USD"
$ws.Range("C30").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.999"
$ws.Range("E30").Value = "  +0.07%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.87"
$ws.Range("E31").Value = "  -1.78%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.05"
$ws.Range("E32").Value = "  +6.34%  "

# Row 33
$ws.Range("E33").Value = "  +7.05%  "

# Row 34
$ws.Range("E34").Value = "  +0.80%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "172.52"
$ws.Range("E35").Value = "  -2.42%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.08"
$ws.Range("E36").Value = "  +12.56%  "

# Row 37
$ws.Range("B37").Value = "PolygonEcosystemToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.405"
$ws.Range("E37").Value = "  +2.13%  "

# Row 38
$ws.Range("B38").Value = "FirstDigitalUSD"
$ws.Range("C38").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.00"
$ws.Range("E38").Value = "  -0.09%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "19.25"
$ws.Range("E39").Value = "  +1.68%  "

# Row 40
$ws.Range("E40").Value = "  +7.61%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "174.79"
$ws.Range("E41").Value = "  +9.67%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.999"
$ws.Range("E42").Value = "  -0.02%  "

# Row 43
$ws.Range("E43").Value = "  +1.10%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "22.14"
$ws.Range("E44").Value = "  +3.26%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0565"
$ws.Range("E45").Value = "  +4.64%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.634"
$ws.Range("E46").Value = "  -0.39%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0965"
$ws.Range("E47").Value = "  +0.18%  "

# Row 48
$ws.Range("E48").Value = "  +1.43%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "18.86"
$ws.Range("E49").Value = "  +4.43%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.74"
$ws.Range("E50").Value = "  +2.88%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "11.32"
$ws.Range("E51").Value = "  -0.73%  "
